# Fix Training Data Issue (#48)
# The "Date" column (BF) was populated with the source filename
# ("5-25-2007-08") instead of the actual game date. NBA stats for a given
# night are reported under the next calendar day, so the correct date is
# one day later than the naive "5-25" read off the file name: 2008-05-25.
#
# Columns BF2:BF31 all hold the same mis-derived text value and need to be
# corrected to the proper ISO date string, written back as literal text
# (not an Excel date serial).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = "5-25-2007-08"
$newValue = "2008-05-25"

$firstRow = 2
$lastRow = 31
$col = 58  # column BF

for ($row = $firstRow; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, $col)
    if ($cell.Value2 -eq $oldValue) {
        # Force the new value to stay a text string instead of letting
        # Excel's input parser reinterpret "2008-05-25" as a date serial,
        # then restore the default "Normal" style so the cell's formatting
        # is unaffected by the temporary text-format nudge.
        $cell.NumberFormat = "@"
        $cell.Value = $newValue
        $cell.Style = "Normal"
    }
}
